# "validacion de detalle de paleta en los kilos netos y brutos"
#
# Row 2 (the first pallet detail row) gets new validated values, and the
# second detail row (row 3) is removed entirely now that only one
# validated row remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (well outside the used range) used to stage values that must
# land in the sheet as literal TEXT. Typing a numeric-looking string (e.g.
# "45") straight into A2/B2 via .Value would get auto-converted to a
# number by Excel's input parser, losing the shared-string/text nature the
# original cells had. Formatting the scratch cell as Text ("@"), entering
# the value there, then Copy / PasteSpecial-Values onto the real target
# cell keeps the text type while leaving the target cell's own style
# (s="0") untouched, since a values-only paste carries no formatting.
$scratchCell = $ws.Range("Z1000")

$scratchCell.NumberFormat = "@"
$scratchCell.Value = "45"
$scratchCell.Copy()
$ws.Range("A2").PasteSpecial(-4163)

$scratchCell.NumberFormat = "@"
$scratchCell.Value = "69"
$scratchCell.Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("C2").Value = "50X50"

# D2's new multi-line text is staged the same way and pasted as values so
# the row-2 auto "fit to wrapped text" row-height recalculation that a
# direct .Value edit on a wrap-text cell would trigger never fires.
$scratchCell.Value = "DFNSA.FSFSDFJÑKÑ.KSAJFH`nSF-SHFKJSHFÑKSADJFSÑKJFH`nSDLFKHSLFKSÑAFHASDÑF"
$scratchCell.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("E2").Value = 850
$ws.Range("F2").Value = 1250

# Drop the scratch row and the now-obsolete third pallet-detail row.
$ws.Rows(1000).Delete()
$ws.Rows(3).Delete()

# Medida / Contenido columns widened slightly to fit the new content.
$ws.Columns(3).ColumnWidth = 7.0
$ws.Columns(4).ColumnWidth = 25.333333333333332
